# Automatic map update — adds "PD" (Q) and "N2" (R) columns to the INCO sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) — mirrors the existing header formatting already
# applied to A1:P1 (bold, thin border, center/top alignment) by copying the
# format from the last existing header cell (P1).
$ws.Range("P1").Copy() | Out-Null
$ws.Range("Q1:R1").PasteSpecial(-4122) | Out-Null
$ws.Range("Q1").Value = "PD"
$ws.Range("R1").Value = "N2"
$excel.CutCopyMode = $false

# Per-row data (Q = PD code, R = N2 / polygon note).
$rows = @(
    @{Row=1; Q='PD'; R='N2'},
    @{Row=2; Q='CON-M'; R='Fuera de Poligono OVL'},
    @{Row=3; Q='BLO-H'; R='Fuera de Poligono OVL'},
    @{Row=4; Q='RET-E'; R='ARATO-25058.PO.1RET'},
    @{Row=5; Q='ATH-Q'; R='Fuera de Poligono OVL'},
    @{Row=6; Q='CON-K'; R='Fuera de Poligono OVL'},
    @{Row=7; Q='AGU-H'; R='Fuera de Poligono OVL'},
    @{Row=8; Q='PUE-?'; R='Fuera de Poligono OVL'},
    @{Row=9; Q='CEN-A'; R='Fuera de Poligono OVL'},
    @{Row=10; Q='CEN-B'; R='Fuera de Poligono OVL'},
    @{Row=11; Q='CEN-G'; R='Fuera de Poligono OVL'},
    @{Row=12; Q='CEN-K'; R='Fuera de Poligono OVL'},
    @{Row=13; Q='NRA-M'; R='Fuera de Poligono OVL'},
    @{Row=14; Q='CEN-J'; R='Fuera de Poligono OVL'},
    @{Row=15; Q='BLO-B'; R='Fuera de Poligono OVL'},
    @{Row=16; Q='CLI-I'; R='Fuera de Poligono OVL'},
    @{Row=17; Q='ALM-A'; R='Fuera de Poligono OVL'},
    @{Row=18; Q='ALM-B'; R='Fuera de Poligono OVL'},
    @{Row=19; Q='ALM-D'; R='Fuera de Poligono OVL'},
    @{Row=20; Q='CEN-P'; R='Fuera de Poligono OVL'},
    @{Row=21; Q='COG-G'; R='Fuera de Poligono OVL'},
    @{Row=22; Q='PUE-J'; R='Fuera de Poligono OVL'},
    @{Row=23; Q='COG-L'; R='Fuera de Poligono OVL'},
    @{Row=24; Q='PUE-K'; R='Fuera de Poligono OVL'},
    @{Row=25; Q='ATH-H'; R='Fuera de Poligono OVL'},
    @{Row=26; Q='BLO-J'; R='Fuera de Poligono OVL'},
    @{Row=27; Q='BLO-I'; R='Fuera de Poligono OVL'},
    @{Row=28; Q='BLO-L'; R='Fuera de Poligono OVL'},
    @{Row=29; Q='PPT-P'; R='Fuera de Poligono OVL'},
    @{Row=30; Q='COG-L'; R='Fuera de Poligono OVL'},
    @{Row=31; Q='BLO-H'; R='Fuera de Poligono OVL'},
    @{Row=32; Q='BLO-H'; R='Fuera de Poligono OVL'},
    @{Row=33; Q='PPT-P'; R='Fuera de Poligono OVL'},
    @{Row=34; Q='PPT-P'; R='Fuera de Poligono OVL'},
    @{Row=35; Q='COG-H'; R='Fuera de Poligono OVL'},
    @{Row=36; Q='PUE-J'; R='Fuera de Poligono OVL'},
    @{Row=37; Q='PPT-K'; R='Fuera de Poligono OVL'},
    @{Row=38; Q='COG-A'; R='Fuera de Poligono OVL'},
    @{Row=39; Q='CON-M'; R='Fuera de Poligono OVL'},
    @{Row=40; Q='CON-K'; R='Fuera de Poligono OVL'}
)

foreach ($item in $rows) {
    if ($item.Row -eq 1) { continue }
    $ws.Cells.Item($item.Row, 17).Value = $item.Q
    $ws.Cells.Item($item.Row, 18).Value = $item.R
}
